# Update mods data [2026-01-14 15:14:53]
#
# Appends a new row (row 65) to the ModCounts sheet with the latest mod
# count reading, mirroring the pattern of all the preceding daily rows:
#   A65 = "2026/01/14"   (text date, same shape as the other Date cells)
#   B65 = "逃离鸭科夫"      (game name, unchanged from every other row)
#   C65 = 1144            (numeric mod count)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Copy formatting only from the previous row (row 64) onto the new
#        row (row 65) so the new cells get the exact same style (centered
#        alignment etc.) as the rest of the table. -4122 = xlPasteFormats.
$ws.Range("A64:C64").Copy()
$ws.Range("A65:C65").PasteSpecial(-4122)

# --- 2. Fill in the Date cell (A65).
#        Assigning a date-shaped string straight to .Value makes Excel
#        "helpfully" reinterpret it as a real date serial number, which is
#        not what the source data does (every Date cell in this sheet is
#        plain text, e.g. "2026/01/13"). To keep it as literal text we
#        build the string with a formula in a scratch cell far away from
#        the real data, copy that computed value, and paste-special just
#        the value (-4163 = xlPasteValues) into A65 - this bypasses the
#        automatic text->date coercion that a direct .Value assignment
#        would trigger. The scratch cell is cleared again afterwards so it
#        leaves no trace in the saved workbook.
$scratch = $ws.Cells.Item(1000, 1)
$scratch.Formula = "=""2026/01/14"""
$scratch.Copy()
$ws.Range("A65").PasteSpecial(-4163)
$scratch.Value = ""

# --- 3. Fill in the remaining cells normally.
$ws.Cells.Item(65, 2).Value = "逃离鸭科夫"
$ws.Cells.Item(65, 3).Value = 1144
